$wb = $excel.ActiveWorkbook

$wsBacklog = $wb.Worksheets.Item("product backlog")
$wsSprint = $wb.Worksheets.Item("sprint backlog")

# Add new user story rows (11-13) to the product backlog sheet
$wsBacklog.Range("E11").Value = "to see pdf reports for employees per department etc."

$wsBacklog.Range("E12").Value = "setup virtual environment for the app"
$wsBacklog.Range("F12").Value = "it can be uploaded and run from servers"

$wsBacklog.Range("E13").Value = "do some performance profiling for the application"
$wsBacklog.Range("F13").Value = "it can be optimized"
$wsBacklog.Range("G13").Value = "http://docs.sqlalchemy.org/en/rel_1_0/faq/performance.html#faq-how-to-profile"

$wsBacklog.Range("A11:I13").RowHeight = 30

# The sprint backlog sheet should no longer be the tab shown/selected,
# keep its existing selection at G5
$wsSprint.Range("G5").Select()

# Make the product backlog sheet the active/selected tab with G13 selected
$wsBacklog.Activate()
$wsBacklog.Range("G13").Select()
